# Auto update: 2025-11-29 03:23:53
# Applies updated market/decision metrics to Sheet1 rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TSM)
$ws.Range("D2").Value = 290.89
$ws.Range("F2").Value = 4.83
$ws.Range("K2").Value = 64.2
$ws.Range("N2").Value = 85.96878041621773

# Row 3 (ASML)
$ws.Range("D3").Value = 1060
$ws.Range("F3").Value = 8.050000000000001
$ws.Range("K3").Value = 54.8
$ws.Range("N3").Value = 85.96878041621773

# Row 4 (AMD)
$ws.Range("D4").Value = 217.53
$ws.Range("F4").Value = 5.59
$ws.Range("K4").Value = 53
$ws.Range("N4").Value = 85.96878041621773

# Row 5 (NVDA)
$ws.Range("D5").Value = 177
$ws.Range("F5").Value = -2.02
$ws.Range("K5").Value = 52.8
$ws.Range("N5").Value = 85.96878041621773

# Row 6 (QCOM)
$ws.Range("D6").Value = 168.09
$ws.Range("F6").Value = 5.33
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 43
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 49
$ws.Range("N6").Value = 85.96878041621773
